# Apply the "TEST101..TEST105" additions + red-flag restyle to the DB sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# --- Add 5 new student codes below the existing list (rows 52-56) ---
$ws.Cells.Item(52,1).Value = "TEST101"
$ws.Cells.Item(53,1).Value = "TEST102"
$ws.Cells.Item(54,1).Value = "TEST103"
$ws.Cells.Item(55,1).Value = "TEST104"
$ws.Cells.Item(56,1).Value = "TEST105"

# --- Restyle rows 49-50 ("MRTiohn","MRTI019") to the red/"done" look    ---
# --- already used by rows 51+ (style index 2), and apply the same red  ---
# --- style to the first new row (52), matching the rest (53-56) with   ---
# --- the plain "not yet" style already used by rows 28-48 (style 3).   ---
# Copy formatting (not values) from cells that already carry the target
# style so we reuse the existing style entries instead of minting new
# font/xf records.
$ws.Cells.Item(51,1).Copy()
$ws.Cells.Item(49,1).PasteSpecial(-4122)
$ws.Cells.Item(50,1).PasteSpecial(-4122)
$ws.Cells.Item(52,1).PasteSpecial(-4122)

$ws.Cells.Item(48,1).Copy()
$ws.Cells.Item(53,1).PasteSpecial(-4122)
$ws.Cells.Item(54,1).PasteSpecial(-4122)
$ws.Cells.Item(55,1).PasteSpecial(-4122)
$ws.Cells.Item(56,1).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Reflect the updated scroll position / selection on the DB sheet ---
$ws.Activate()
$ws.Range("D52").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
